$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp values for column C (time column), rows 2..37, in row order.
# These correspond 1:1 to the updated "זמן" (time) values from the re-run
# of the registration validation test suite.
$newTimes = @(
    "01/12/2025 22:34:39",
    "01/12/2025 22:34:43",
    "01/12/2025 22:34:45",
    "01/12/2025 22:34:50",
    "01/12/2025 22:34:55",
    "01/12/2025 22:34:55",
    "01/12/2025 22:35:16",
    "01/12/2025 22:35:25",
    "01/12/2025 22:35:30",
    "01/12/2025 22:35:40",
    "01/12/2025 22:35:44",
    "01/12/2025 22:35:46",
    "01/12/2025 22:35:50",
    "01/12/2025 22:35:57",
    "01/12/2025 22:35:57",
    "01/12/2025 22:36:18",
    "01/12/2025 22:36:21",
    "01/12/2025 22:36:26",
    "01/12/2025 22:36:36",
    "01/12/2025 22:36:40",
    "01/12/2025 22:36:42",
    "01/12/2025 22:36:47",
    "01/12/2025 22:36:52",
    "01/12/2025 22:36:52",
    "01/12/2025 22:37:13",
    "01/12/2025 22:37:16",
    "01/12/2025 22:37:20",
    "01/12/2025 22:37:30",
    "01/12/2025 22:37:34",
    "01/12/2025 22:37:37",
    "01/12/2025 22:37:41",
    "01/12/2025 22:37:52",
    "01/12/2025 22:37:52",
    "01/12/2025 22:38:13",
    "01/12/2025 22:38:16",
    "01/12/2025 22:38:21"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newTimes[$i]
}
